$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($val -eq "H2495\S3") {
            $cell.Value2 = "H2495\S3\EOBO"
        } elseif ($val -eq "H2495\S8") {
            $cell.Value2 = "H2495\S8\EOBO"
        } elseif ($val -eq "H2495\S9") {
            $cell.Value2 = "H2495\S9\EOBN"
        }
    }
}
